$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '36.673.67'
$ws.Cells.Item(2, 5).Value = '  -1.12%  '
$ws.Cells.Item(3, 4).Value = '2.080.32'
$ws.Cells.Item(3, 5).Value = '  +1.59%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '244.35'
$ws.Cells.Item(5, 5).Value = '  -1.36%  '
$ws.Cells.Item(6, 5).Value = '  -2.07%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '53.65'
$ws.Cells.Item(8, 5).Value = '  -5.78%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '58.75'
$ws.Cells.Item(9, 5).Value = '  -2.01%  '
$ws.Cells.Item(10, 5).Value = '  -4.06%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0757'
$ws.Cells.Item(11, 5).Value = '  -2.12%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.109'
$ws.Cells.Item(12, 5).Value = '  +0.86%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '14.91'
$ws.Cells.Item(13, 5).Value = '  -6.05%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.887'
$ws.Cells.Item(14, 5).Value = '  +2.57%  '
$ws.Cells.Item(15, 4).Value = '2.386.19'
$ws.Cells.Item(15, 5).Value = '  +1.58%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '5.46'
$ws.Cells.Item(16, 5).Value = '  -3.33%  '
$ws.Cells.Item(17, 4).Value = '2.044.66'
$ws.Cells.Item(17, 5).Value = '  -0.17%  '
$ws.Cells.Item(18, 4).Value = '36.643.84'
$ws.Cells.Item(18, 5).Value = '  -1.03%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '17.15'
$ws.Cells.Item(19, 5).Value = '  -3.88%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '72.34'
$ws.Cells.Item(20, 5).Value = '  -3.02%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0873'
$ws.Cells.Item(21, 5).Value = '  -1.56%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.40'
$ws.Cells.Item(22, 5).Value = '  +0.91%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '239.13'
$ws.Cells.Item(23, 5).Value = '  +1.06%  '
$ws.Cells.Item(24, 5).Value = '  +0.07%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.38'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '9.77'
$ws.Cells.Item(26, 5).Value = '  +3.85%  '
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '167.49'
$ws.Cells.Item(27, 5).Value = '  -0.72%  '
$ws.Cells.Item(28, 2).Value = 'PancakeSwap'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.14'
$ws.Cells.Item(28, 5).Value = '  -1.18%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '20.57'
$ws.Cells.Item(29, 5).Value = '  +3.10%  '
$ws.Cells.Item(30, 5).Value = '  -1.25%  '
$ws.Cells.Item(31, 5).Value = '  +10.26%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.17'
$ws.Cells.Item(32, 5).Value = '  +4.79%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.63'
$ws.Cells.Item(33, 5).Value = '  +3.85%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0604'
$ws.Cells.Item(34, 5).Value = '  -1.73%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.36'
$ws.Cells.Item(35, 5).Value = '  +5.89%  '
$ws.Cells.Item(36, 5).Value = '  +0.15%  '
$ws.Cells.Item(37, 5).Value = '  +4.37%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0823'
$ws.Cells.Item(38, 5).Value = '  -7.14%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.25'
$ws.Cells.Item(39, 5).Value = '  -5.78%  '
$ws.Cells.Item(40, 5).Value = '  +1.31%  '
$ws.Cells.Item(41, 5).Value = '  -1.01%  '
$ws.Cells.Item(42, 2).Value = 'Cronos'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0953'
$ws.Cells.Item(42, 5).Value = '  +3.07%  '
$ws.Cells.Item(43, 2).Value = 'THORChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '4.77'
$ws.Cells.Item(43, 5).Value = '  -8.25%  '
$ws.Cells.Item(44, 2).Value = 'HuobiToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.85'
$ws.Cells.Item(44, 5).Value = '  -10.19%  '
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '95.63'
$ws.Cells.Item(45, 5).Value = '  +0.20%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '15.88'
$ws.Cells.Item(46, 5).Value = '  -7.18%  '
$ws.Cells.Item(47, 4).Value = '1.370.74'
$ws.Cells.Item(47, 5).Value = '  +8.11%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '7.27'
$ws.Cells.Item(48, 5).Value = '  +7.08%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.41'
$ws.Cells.Item(49, 5).Value = '  +0.04%  '
$ws.Cells.Item(50, 5).Value = '  +1.53%  '
$ws.Cells.Item(51, 4).Value = '2.269.45'
$ws.Cells.Item(51, 5).Value = '  +1.63%  '
